$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the 7-day block starting at row 141 forward by one day ---
# B141 holds the literal start date; C141:H141 are formulas (=prev+1) that
# recalculate automatically once B141 changes.
$ws.Range("B141").Value = $ws.Range("B141").Value2 + 1

# --- New footnote question added below the table first ---
$ws.Range("B156").Value = "未在床上/卧室使用电子设备计算入吗？"

# --- Fill in day-1 answers for the newly-current week (rows 142-155) ---
$ws.Range("B142").Value = "7：42"
$ws.Range("B143").Value = "7：42"
$ws.Range("B144").Value = "23：00"
$ws.Range("B145").Value = "23：20"
$ws.Range("B146").Value = 0
$ws.Range("B147").Value = 2
$ws.Range("B148").Value = 30
$ws.Range("B149").Value = 480
$ws.Range("B150").Value = "无"
$ws.Range("B151").Value = 0
$ws.Range("B152").Value = 3
$ws.Range("B153").Value = 4
$ws.Range("B154").Value = 4
$ws.Range("B155").Value = "有 50 min"

# --- Keep the view in sync with where the user ended up editing ---
$ws.Range("B155").Select()
